$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("info")

# Widen column B (dimension grows to A1:D8 automatically once row 8 is filled)
$ws.Columns.Item(2).ColumnWidth = 41.285714285714285

# New row 8: dbdriver / driver class name
$ws.Range("A8").Value = "dbdriver"
$ws.Range("B8").Value = "com.microsoft.sqlserver.jdbc.SQLServerDriver"

# Update the active selection shown when the sheet is (re)opened
$ws.Range("B3").Select() | Out-Null
